# Apply the edits described by the diff for CasosColombia.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) I14: was the shared string "NaN" -> becomes the number 1
$ws.Range("I14").Value = 1

# 2) AI185: was the number 308 -> becomes the shared string "NaN"
$ws.Range("AI185").Value = "NaN"

# 3) Append a brand-new data row 204 with values for columns A..DX
$rowValues = @(
    44098, 790823, 2741, 107933, 66938, 259008, 28374, 6453, 5239, 8036, 8790, 18129, 3975,
    23366, 31877, 7689, 10204, 14864, 14138, 17825, 15001, 3695, 2948, 10166, 29060, 13867,
    11506, 59212, 2012, 1079, 733, 471, 736, 465, 657, 2042, 5373, 37792, 9473,
    2543, 45984, 1100, 22693, 1524, 10186, 1657, 1603, 7954, 2001, 962, 2500, 2677,
    62599, 13947, 6050, 9644, 6835, 257, 1464, 2723, 743, 2153, 9718, 9515, 10461,
    14272, 1964, 899, 13237, 10939, 12770, 2846, 2137, 5624, 4714, 2137, 5769, 3575,
    2126, 968, 2946, 2225, 1942, 1726, 6250, 2142, 1447, 1777, 2095, 2161, 2543,
    1667, 1213, 1209, 976, 3392, 1456, 949, 1061, 1718, 1581, 806, 887, 1304,
    1633, 1511, 1568, 1226, 334, 368, 819, 766, 490, 543, 381, 669, 751,
    527, 492, 374, 521, 135994, 334907, 18626, 145156, 89802, 44110, 12467
)

$targetRow = 204
$col = 1
foreach ($v in $rowValues) {
    $ws.Cells.Item($targetRow, $col).Value = $v
    $col = $col + 1
}

# 4) Update the view so the active/selected cell becomes B204, matching
#    the workbook as saved after the new row was entered.
$ws.Range("B204").Select()
